$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '91.414.26'
$ws.Range('E2').Value = '  +0.99%  '
$ws.Range('D3').Value = '3.169.95'
$ws.Range('E3').Value = '  +1.83%  '
$ws.Range('E4').Value = '  +0.24%  '
$ws.Range('D5').Value = '''239.14'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +2.16%  '
$ws.Range('D6').Value = '''620.74'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.45%  '
$ws.Range('D7').Value = '''1.14'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +6.82%  '
$ws.Range('D8').Value = '''0.372'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +1.26%  '
$ws.Range('E9').Value = '  +0.01%  '
$ws.Range('D10').Value = '3.167.57'
$ws.Range('E10').Value = '  +1.95%  '
$ws.Range('D11').Value = '''0.744'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.88%  '
$ws.Range('D12').Value = '''0.206'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +3.90%  '
$ws.Range('D13').Value = '''0.0000248'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -1.40%  '
$ws.Range('D14').Value = '''35.40'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -1.68%  '
$ws.Range('E15').Value = '  +0.99%  '
$ws.Range('D16').Value = '91.485.82'
$ws.Range('E16').Value = '  +1.37%  '
$ws.Range('D17').Value = '3.748.12'
$ws.Range('D18').Value = '3.137.34'
$ws.Range('E18').Value = '  +1.35%  '
$ws.Range('E19').Value = '  -4.25%  '
$ws.Range('D20').Value = '''15.34'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +9.16%  '
$ws.Range('D21').Value = '''5.91'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +6.12%  '
$ws.Range('D22').Value = '''0.0000210'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -5.90%  '
$ws.Range('D23').Value = '''444.16'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.49%  '
$ws.Range('D24').Value = '''9.22'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +2.96%  '
$ws.Range('D25').Value = '''6.10'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +3.08%  '
$ws.Range('D26').Value = '''89.07'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.37%  '
$ws.Range('D27').Value = '''12.11'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.15%  '
$ws.Range('D28').Value = '3.349.26'
$ws.Range('E28').Value = '  +3.20%  '
$ws.Range('E29').Value = '  -0.16%  '
$ws.Range('D30').Value = '''0.131'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +45.17%  '
$ws.Range('D31').Value = '''0.171'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +6.88%  '
$ws.Range('D32').Value = '''0.228'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +17.09%  '
$ws.Range('D33').Value = '''9.38'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.45%  '
$ws.Range('E34').Value = '  +9.57%  '
$ws.Range('D35').Value = '''7.74'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +6.08%  '
$ws.Range('D36').Value = '''26.50'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +2.40%  '
$ws.Range('D37').Value = '''0.887'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -11.20%  '
$ws.Range('D38').Value = '''510.50'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +1.22%  '
$ws.Range('B39').Value = 'Fetch.AI'
$ws.Range('C39').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D39').Value = '''1.36'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +5.81%  '
$ws.Range('B40').Value = 'PancakeSwap'
$ws.Range('C40').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D40').Value = '''1.94'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +1.64%  '
$ws.Range('D41').Value = '''0.460'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +13.14%  '
$ws.Range('D42').Value = '''3.81'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +10.28%  '
$ws.Range('D43').Value = '''3.44'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -10.39%  '
$ws.Range('D44').Value = '''22.15'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.13%  '
$ws.Range('D46').Value = '''0.718'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +3.72%  '
$ws.Range('D47').Value = '''156.65'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +2.65%  '
$ws.Range('D48').Value = '''1.92'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.91%  '
$ws.Range('D49').Value = '''1.39'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +4.00%  '
$ws.Range('D50').Value = '''0.0328'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +12.25%  '
$ws.Range('D51').Value = '''4.46'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.72%  '
